$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# --- Rows whose "D" (Price) value looks like a plain single-dot number and
#     must be forced back to text so Excel does not coerce it into a Double ---
Set-TextValue "D5"  "213.71"
Set-TextValue "D8"  "24.13"
Set-TextValue "D11" "0.0890"
Set-TextValue "D15" "3.75"
Set-TextValue "D17" "63.15"
Set-TextValue "D18" "227.85"
Set-TextValue "D20" "7.48"
Set-TextValue "D23" "9.32"
Set-TextValue "D25" "151.83"
Set-TextValue "D26" "15.19"
Set-TextValue "D28" "6.59"
Set-TextValue "D33" "3.14"
Set-TextValue "D37" "2.34"
Set-TextValue "D40" "0.540"
Set-TextValue "D41" "0.814"
Set-TextValue "D45" "0.989"
Set-TextValue "D46" "64.24"

# --- Rows whose "D" value already contains two dots (or non-numeric
#     characters) so plain assignment keeps it as text ---
$ws.Range("D2").Value = "28.253.35"
$ws.Range("D3").Value = "1.590.71"
$ws.Range("D12").Value = "1.817.91"
$ws.Range("D13").Value = "1.588.09"
$ws.Range("D16").Value = "28.302.26"
$ws.Range("D19").Value = "0.0₃0710"
$ws.Range("D34").Value = "1.399.13"
$ws.Range("D47").Value = "1.730.83"

# --- Volume(1h) column updates (E) - always safe text, never pure numbers ---
$ws.Range("E2").Value  = "  +3.66%  "
$ws.Range("E3").Value  = "  +1.72%  "
$ws.Range("E4").Value  = "  +0.13%  "
$ws.Range("E5").Value  = "  +1.27%  "
$ws.Range("E6").Value  = "  +0.81%  "
$ws.Range("E7").Value  = "  -0.02%  "
$ws.Range("E8").Value  = "  +8.73%  "
$ws.Range("E9").Value  = "  +0.49%  "
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  +3.82%  "
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("E18").Value = "  +4.49%  "
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").Value = "  -4.00%  "
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("E36").Value = "  -7.12%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  +8.56%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("E44").Value = "  +6.74%  "
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("E47").Value = "  +1.76%  "

# --- Rows 48-51: two pairs of coins swap rank order, with refreshed values ---
# Row 48 becomes Quant (previously row 49)
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D48" "87.60"
$ws.Range("E48").Value = "  +2.02%  "

# Row 49 becomes mCoin (previously row 48)
$ws.Range("B49").Value = "mCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
Set-TextValue "D49" "2.14"
$ws.Range("E49").Value = "  +1.23%  "

# Row 50 becomes BabyDogeCoin (previously row 51)
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0103"
$ws.Range("E50").Value = "  -1.50%  "

# Row 51 becomes Cronos (previously row 50)
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.0524"
$ws.Range("E51").Value = "  +0.05%  "
